# Scheduled-runner update: refresh market-board price/profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2799.763
$ws.Range("I15").Value = 2799.763
$ws.Range("K15").Value = 8399.289000000001
$ws.Range("M15").Value = -8230.289000000001

$ws.Range("H40").Value = 2564.8572
$ws.Range("I40").Value = 3461.2
$ws.Range("J40").Value = 1750
$ws.Range("K40").Value = 3461.2
$ws.Range("L40").Value = 1750
$ws.Range("M40").Value = -3286.2
$ws.Range("N40").Value = -2100

$ws.Range("H55").Value = 216.3077
$ws.Range("I55").Value = 111.833336
$ws.Range("J55").Value = 305.85715
$ws.Range("K55").Value = 111.833336
$ws.Range("L55").Value = 305.85715
$ws.Range("M55").Value = 102.166664
$ws.Range("N55").Value = -733.85715

$ws.Range("H64").Value = 2833.1143
$ws.Range("I64").Value = 3078.625
$ws.Range("J64").Value = 2760.3704
$ws.Range("K64").Value = 3078.625
$ws.Range("L64").Value = 2760.3704
$ws.Range("M64").Value = -2830.625
$ws.Range("N64").Value = -3256.3704

$ws.Range("H67").Value = 2833.1143
$ws.Range("I67").Value = 3078.625
$ws.Range("J67").Value = 2760.3704
$ws.Range("K67").Value = 3078.625
$ws.Range("L67").Value = 2760.3704
$ws.Range("M67").Value = -2220.625
$ws.Range("N67").Value = -4476.3704

$ws.Range("H74").Value = 3029.4194
$ws.Range("I74").Value = 3221.1333
$ws.Range("J74").Value = 2849.6875
$ws.Range("K74").Value = 3221.1333
$ws.Range("L74").Value = 2849.6875
$ws.Range("M74").Value = -2285.1333
$ws.Range("N74").Value = -4721.6875

$ws.Range("H77").Value = 3029.4194
$ws.Range("I77").Value = 3221.1333
$ws.Range("J77").Value = 2849.6875
$ws.Range("K77").Value = 16105.6665
$ws.Range("L77").Value = 14248.4375
$ws.Range("M77").Value = -11425.6665
$ws.Range("N77").Value = -23608.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2826.7727
$ws.Range("I63").Value = 2020.7142
$ws.Range("K63").Value = 2020.7142
$ws.Range("M63").Value = -1334.7142

$ws.Range("H66").Value = 2826.7727
$ws.Range("I66").Value = 2020.7142
$ws.Range("K66").Value = 10103.571
$ws.Range("M66").Value = -6671.571

$ws.Range("H122").Value = 2059.8
$ws.Range("I122").Value = 1088.75
$ws.Range("J122").Value = 5944
$ws.Range("K122").Value = 3266.25
$ws.Range("L122").Value = 17832
$ws.Range("M122").Value = -816.25
$ws.Range("N122").Value = -22732

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1608.6364
$ws.Range("I20").Value = 1402.3334
$ws.Range("J20").Value = 1780.5555
$ws.Range("K20").Value = 1402.3334
$ws.Range("L20").Value = 1780.5555
$ws.Range("M20").Value = -1155.3334
$ws.Range("N20").Value = -2274.5555

$ws.Range("H22").Value = 160
$ws.Range("I22").Value = 160
$ws.Range("K22").Value = 160
$ws.Range("M22").Value = 13

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 658
$ws.Range("I22").Value = 496.66666
$ws.Range("K22").Value = 496.66666
$ws.Range("M22").Value = -146.66666

$ws.Range("H33").Value = 2010.3334
$ws.Range("I33").Value = 1515.5
$ws.Range("J33").Value = 3000
$ws.Range("K33").Value = 1515.5
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = -1136.5
$ws.Range("N33").Value = -3758

$ws.Range("H44").Value = 12000
$ws.Range("I44").Value = 12000
$ws.Range("K44").Value = 12000
$ws.Range("M44").Value = -11558

$ws.Range("H107").Value = 41668850
$ws.Range("I107").Value = 83335700
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 83335700
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -83333780
$ws.Range("N107").Value = -5840

$ws.Range("H134").Value = 4593.826
$ws.Range("I134").Value = 4576.737
$ws.Range("J134").Value = 4675
$ws.Range("K134").Value = 13730.211
$ws.Range("L134").Value = 14025
$ws.Range("M134").Value = -11195.211
$ws.Range("N134").Value = -19095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1905863.1
$ws.Range("I131").Value = 22222508
$ws.Range("J131").Value = 1177.6562
$ws.Range("K131").Value = 66667524
$ws.Range("L131").Value = 3532.9686
$ws.Range("M131").Value = -66662484
$ws.Range("N131").Value = -13612.9686

$ws.Range("H139").Value = 5197.0625
$ws.Range("I139").Value = 3788.25
$ws.Range("J139").Value = 5666.6665
$ws.Range("K139").Value = 11364.75
$ws.Range("L139").Value = 16999.9995
$ws.Range("M139").Value = -6224.75
$ws.Range("N139").Value = -27279.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2508.2144
$ws.Range("I102").Value = 2384.5833
$ws.Range("J102").Value = 3250
$ws.Range("K102").Value = 2384.5833
$ws.Range("L102").Value = 3250
$ws.Range("M102").Value = -762.5832999999998
$ws.Range("N102").Value = -6494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3846894.5
$ws.Range("I22").Value = 9615708
$ws.Range("J22").Value = 1018.75
$ws.Range("K22").Value = 9615708
$ws.Range("L22").Value = 1018.75
$ws.Range("M22").Value = -9615413
$ws.Range("N22").Value = -1608.75

$ws.Range("H27").Value = 3846894.5
$ws.Range("I27").Value = 9615708
$ws.Range("J27").Value = 1018.75
$ws.Range("K27").Value = 9615708
$ws.Range("L27").Value = 1018.75
$ws.Range("M27").Value = -9615601
$ws.Range("N27").Value = -1232.75

$ws.Range("H33").Value = 2500
$ws.Range("I33").Value = 2000
$ws.Range("J33").Value = 3000
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = -1710
$ws.Range("N33").Value = -3580

$ws.Range("H55").Value = 7936719
$ws.Range("I55").Value = 12987209
$ws.Range("J55").Value = 234.42857
$ws.Range("K55").Value = 12987209
$ws.Range("L55").Value = 234.42857
$ws.Range("M55").Value = -12987036
$ws.Range("N55").Value = -580.42857

$ws.Range("H57").Value = 3041
$ws.Range("I57").Value = 3041
$ws.Range("K57").Value = 3041
$ws.Range("M57").Value = -2475

$ws.Range("H122").Value = 2592.261
$ws.Range("I122").Value = 2408.0588
$ws.Range("J122").Value = 3114.1667
$ws.Range("K122").Value = 7224.176399999999
$ws.Range("L122").Value = 9342.500100000001
$ws.Range("M122").Value = -4774.176399999999
$ws.Range("N122").Value = -14242.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H122").Value = 985.1875
$ws.Range("I122").Value = 705.3333
$ws.Range("J122").Value = 1824.75
$ws.Range("K122").Value = 2115.9999
$ws.Range("L122").Value = 5474.25
$ws.Range("M122").Value = 334.0001000000002
$ws.Range("N122").Value = -10374.25
